$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.094.19"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").Value = "2.050.72"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'231.36"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("D6").Value = "'0.616"
$ws.Range("E6").Value = "  +2.71%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'57.20"
$ws.Range("E8").Value = "  +3.16%  "
$ws.Range("D9").Value = "'0.380"
$ws.Range("E9").Value = "  +2.14%  "
$ws.Range("D10").Value = "'57.80"
$ws.Range("E10").Value = "  +1.06%  "
$ws.Range("D11").Value = "'0.0754"
$ws.Range("E11").Value = "  +0.94%  "
$ws.Range("E12").Value = "  +0.83%  "
$ws.Range("D13").Value = "2.348.34"
$ws.Range("E13").Value = "  +0.87%  "
$ws.Range("D14").Value = "'14.25"
$ws.Range("E14").Value = "  -0.77%  "
$ws.Range("D15").Value = "'20.73"
$ws.Range("E15").Value = "  +2.26%  "
$ws.Range("D16").Value = "'0.770"
$ws.Range("E16").Value = "  +0.73%  "
$ws.Range("D17").Value = "'5.13"
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").Value = "2.047.99"
$ws.Range("E18").Value = "  +1.03%  "
$ws.Range("D19").Value = "37.015.08"
$ws.Range("E19").Value = "  +0.67%  "
$ws.Range("D20").Value = "'6.28"
$ws.Range("E20").Value = "  +13.52%  "
$ws.Range("D21").Value = "'68.79"
$ws.Range("E21").Value = "  +1.83%  "
$ws.Range("D22").Value = "0.0₃0806"
$ws.Range("E22").Value = "  +0.91%  "
$ws.Range("D23").Value = "'224.11"
$ws.Range("E23").Value = "  +1.39%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").Value = "'2.43"
$ws.Range("E25").Value = "  +1.71%  "
$ws.Range("E26").Value = "  -0.52%  "
$ws.Range("D27").Value = "'165.25"
$ws.Range("E27").Value = "  +1.45%  "
$ws.Range("D28").Value = "'1.45"
$ws.Range("E28").Value = "  +7.71%  "
$ws.Range("D29").Value = "'8.74"
$ws.Range("E29").Value = "  +0.65%  "
$ws.Range("D30").Value = "'18.96"
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("D31").Value = "'0.125"
$ws.Range("E31").Value = "  -3.18%  "
$ws.Range("D32").Value = "'0.117"
$ws.Range("E32").Value = "  -0.37%  "
$ws.Range("D33").Value = "'4.43"
$ws.Range("E33").Value = "  +1.02%  "
$ws.Range("D36").Value = "'4.52"
$ws.Range("E36").Value = "  +5.74%  "
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("E38").Value = "  -0.78%  "
$ws.Range("D39").Value = "'3.25"
$ws.Range("E39").Value = "  -1.24%  "
$ws.Range("E40").Value = "  -2.71%  "
$ws.Range("D41").Value = "'4.50"
$ws.Range("E41").Value = "  +8.41%  "
$ws.Range("E42").Value = "  +1.24%  "
$ws.Range("D43").Value = "1.480.99"
$ws.Range("D44").Value = "'95.89"
$ws.Range("E44").Value = "  +3.15%  "
$ws.Range("E45").Value = "  -0.71%  "
$ws.Range("E46").Value = "  +2.97%  "
$ws.Range("E47").Value = "  +2.66%  "
$ws.Range("D50").Value = "'7.10"
$ws.Range("E50").Value = "  +3.00%  "
$ws.Range("D51").Value = "'2.94"
$ws.Range("E51").Value = "  +1.11%  "

# Row 34 <-> Row 35 content swap (Hedera / LidoDAOToken)
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "'2.53"
$ws.Range("E34").Value = "  +2.47%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.0610"
$ws.Range("E35").Value = "  +1.08%  "

# Row 48 <-> Row 49 content swap (ARBITRUM / InjectiveProtocol)
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'15.20"
$ws.Range("E48").Value = "  -3.15%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "'1.02"
$ws.Range("E49").Value = "  +1.04%  "
